$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("testCitizen")

# --- testCitizen (Sheet2): columns A and B swapped + renamed ---
$ws2.Range("A1").Value = "ulke15755"
$ws2.Range("B1").Value = "u1535"
$ws2.Range("A2").Value = "ulke15756"
$ws2.Range("B2").Value = "u1536"
$ws2.Range("A3").Value = "ulke15757"
$ws2.Range("B3").Value = "u1537"
$ws2.Range("A4").Value = "ulke15758"
$ws2.Range("B4").Value = "u1538"
$ws2.Range("A5").Value = "ulke15759"
$ws2.Range("B5").Value = "u1539"
$ws2.Range("A6").Value = "ulke15760"
$ws2.Range("B6").Value = "u1540"
$ws2.Range("A7").Value = "ulke15761"
$ws2.Range("B7").Value = "u1541"
$ws2.Range("A8").Value = "ulke15762"
$ws2.Range("B8").Value = "u1542"

# --- view/selection changes ---
# Sheet1 selection -> G20, while keeping testCitizen as the active/tab-selected sheet
$ws1.Range("G20").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("A1:A8").Select() | Out-Null

# Zoom on testCitizen
$excel.ActiveWindow.Zoom = 145
